$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 180, shifting existing rows 180-199 down to 181-200.
$ws.Rows(180).Insert()

# Populate the newly inserted row 180 with the new record.
# Columns A,B,C,E,F,G,H,I,R are identical across every row of this sheet
# (same market/product), so copy them from the row now sitting at 181.
# NOTE: read cell contents via the Value() method call (not the bare
# .Value property) - the property accessor surfaces a COM reflection
# descriptor string here instead of the actual cell value.
$ws.Cells.Item(180, 1).Value = $ws.Cells.Item(181, 1).Value()    # A: Mercado ID
$ws.Cells.Item(180, 2).Value = $ws.Cells.Item(181, 2).Value()    # B: Mercado
$ws.Cells.Item(180, 3).Value = $ws.Cells.Item(181, 3).Value()    # C: Región

$ws.Range("D180").Value = 44776                                  # D: Fecha

$ws.Cells.Item(180, 5).Value = $ws.Cells.Item(181, 5).Value()    # E: Codreg
$ws.Cells.Item(180, 6).Value = $ws.Cells.Item(181, 6).Value()    # F: Categoría ID
$ws.Cells.Item(180, 7).Value = $ws.Cells.Item(181, 7).Value()    # G: Categoría
$ws.Cells.Item(180, 8).Value = $ws.Cells.Item(181, 8).Value()    # H: Variedad
$ws.Cells.Item(180, 9).Value = $ws.Cells.Item(181, 9).Value()    # I: Calidad

$ws.Range("J180").Value = 100                                    # J: Volumen
$ws.Range("K180").Value = 1200                                   # K: Precio mínimo
$ws.Range("L180").Value = 1500                                   # L: Precio máximo
$ws.Range("M180").Value = 1350                                   # M: Precio promedio ponderado
$ws.Range("N180").Value = "$/atado 1 a 1,5 kilos"                 # N: Unidad de comercialización
$ws.Range("O180").Value = "Región de Los Lagos"                   # O: Origen
$ws.Range("P180").Value = 900                                    # P: Precio $/Kg
$ws.Range("Q180").Value = 1.5                                    # Q: Kg o Unidades

$ws.Cells.Item(180, 18).Value = $ws.Cells.Item(181, 18).Value()  # R: Hortaliza/Fruta
